$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ H = "19.440.696/0001-03"; I = "Não Disponível"; J = "Logradouro:Rua Fabio Modesto, 158"; K = "Vila Joaquim Sales" }
    3  = @{ H = "09.371.163/0001-55"; I = "LABOMINAS"; J = "Logradouro:Rua A, Alphaville Ii"; K = "Coqueiro" }
    4  = @{ H = "09.371.163/0001-55"; I = "LABOMINAS"; J = "Logradouro:Rua A, Alphaville Ii"; K = "Coqueiro" }
    5  = @{ H = "05.642.544/0001-70"; I = "Não Disponível"; J = "Logradouro:Rua Uruguai, 533"; K = "Centro"; L = "85.805-010" }
    6  = @{ H = "05.642.544/0001-70"; I = "Não Disponível"; J = "Logradouro:Rua Uruguai, 533"; K = "Centro"; L = "85.805-010" }
    8  = @{ H = "13.280.276/0001-12" }
    9  = @{ H = "46.139.952/0001-91"; J = "Logradouro:Rua Padre Joao, 11-25"; K = "Vila Santa Tereza"; L = "17.012-020" }
    10 = @{ H = "46.139.952/0001-91"; J = "Logradouro:Rua Padre Joao, 11-25"; K = "Vila Santa Tereza"; L = "17.012-020" }
    11 = @{ H = "61.486.650/0001-83"; I = "Não Disponível"; J = "Logradouro:Avenida Jurua, 434"; K = "Alphaville Centro Industrial e Empresarial/alphaville."; L = "06.455-010" }
    12 = @{ H = "60.160.546/0001-31"; I = "DIGIMED"; J = "Logradouro:Rua dos Marianos, 227"; K = "Vila Gea"; L = "04.691-110" }
    13 = @{ H = "60.160.546/0001-31"; I = "DIGIMED"; J = "Logradouro:Rua dos Marianos, 227"; K = "Vila Gea"; L = "04.691-110" }
    14 = @{ H = "57.442.774/0001-90"; I = "Não Disponível"; J = "Logradouro:Rua Pedro Stancato, 860"; K = "Chacaras Campos dos Amarais"; L = "13.082-050" }
    16 = @{ H = "33.435.231/0001-87"; I = "GE CELMA"; J = "Logradouro:Rua Alice Herve, 356"; K = "Bingen"; L = "25.665-010" }
    17 = @{ H = "33.435.231/0001-87"; I = "GE CELMA"; J = "Logradouro:Rua Alice Herve, 356"; K = "Bingen"; L = "25.665-010" }
    18 = @{ H = "33.435.231/0001-87"; I = "GE CELMA"; J = "Logradouro:Rua Alice Herve, 356"; K = "Bingen"; L = "25.665-010" }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
